# "Updated photodetection in FD and TD"
# Slide 4 ("Progress"): bump the title font size and the bulleted
# content font sizes. PowerPoint auto-adds <a:normAutofit/> to the
# text body once the (now larger) text needs to shrink-to-fit, which
# happens automatically as a side effect of the Font.Size assignments
# below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Title placeholder: "Progress" -> 44pt ---------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Font.Size = 44

# --- Content placeholder: per-bullet font sizes -----------------------
$content = $s.Shapes.Item(3)
$body = $content.TextFrame.TextRange

# Top-level bullets (28pt)
$body.Paragraphs(1).Font.Size = 28   # Photodetection PCB design and implementation
$body.Paragraphs(2).Font.Size = 28   # TTL generation script
$body.Paragraphs(3).Font.Size = 28   # Sensing setup functionality testing

# Second-level (indented) bullets (24pt)
$body.Paragraphs(4).Font.Size = 24   # Driver tests
$body.Paragraphs(5).Font.Size = 24   # Lock-in tests
